$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.133.41"
$ws.Range("E2").Value = "  -0.90%  "

$ws.Range("D3").Value = "3.577.96"
$ws.Range("E3").Value = "  -1.67%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.31"
$ws.Range("E5").Value = "  -2.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.57"
$ws.Range("E6").Value = "  -4.41%  "

$ws.Range("D7").Value = "3.572.22"
$ws.Range("E7").Value = "  -1.62%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.620"
$ws.Range("E8").Value = "  -3.96%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("E10").Value = "  -1.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.652"
$ws.Range("E11").Value = "  -4.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.20"
$ws.Range("E12").Value = "  -5.06%  "

$ws.Range("E13").Value = "  -3.60%  "

$ws.Range("E14").Value = "  -4.43%  "

$ws.Range("D15").Value = "4.150.88"
$ws.Range("E15").Value = "  -1.74%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.68"
$ws.Range("E16").Value = "  -4.34%  "

$ws.Range("D17").Value = "3.580.51"
$ws.Range("E17").Value = "  -1.59%  "

$ws.Range("D18").Value = "70.062.32"
$ws.Range("E18").Value = "  -1.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.55"
$ws.Range("E19").Value = "  -1.94%  "

$ws.Range("E20").Value = "  -0.98%  "

$ws.Range("E21").Value = "  -3.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "491.91"
$ws.Range("E22").Value = "  +0.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.49"
$ws.Range("E23").Value = "  +0.82%  "

$ws.Range("E24").Value = "  -5.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.82"
$ws.Range("E25").Value = "  +5.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.40"
$ws.Range("E26").Value = "  -2.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.45"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.97"
$ws.Range("E28").Value = "  -6.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.36"
$ws.Range("E29").Value = "  -2.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.76"
$ws.Range("E30").Value = "  -2.17%  "

$ws.Range("E31").Value = "  -3.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.13"
$ws.Range("E32").Value = "  -1.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "65.86"
$ws.Range("E33").Value = "  -0.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.116"
$ws.Range("E34").Value = "  -5.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "576.97"
$ws.Range("E35").Value = "  -5.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.21"
$ws.Range("E36").Value = "  +13.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "39.02"
$ws.Range("E37").Value = "  -3.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.409"
$ws.Range("E38").Value = "  -0.99%  "

$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("D40").Value = "0.0₃0795"
$ws.Range("E40").Value = "  -4.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.50"
$ws.Range("E41").Value = "  -1.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.16"
$ws.Range("E42").Value = "  -0.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.134"
$ws.Range("E43").Value = "  -9.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.64"
$ws.Range("E44").Value = "  +8.73%  "

$ws.Range("E45").Value = "  -4.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0456"
$ws.Range("E46").Value = "  -0.96%  "

$ws.Range("D47").Value = "3.182.88"
$ws.Range("E47").Value = "  -4.37%  "

$ws.Range("E48").Value = "  -2.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.59"
$ws.Range("E49").Value = "  +32.23%  "

$ws.Range("E50").Value = "  -2.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  -0.09%  "
